# feat : make ExerciseRecords_api and related form
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "exercise" API row (row 12), which was previously blank.
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "POST"
$ws.Range("D12").Value = "운동 데이터 추가"
$ws.Range("C12").Value = "/api6/create_exercise/"

# Update the active selection to reflect where the author ended up working.
$ws.Range("C13").Select()
